# Update "paises" (countries) COVID data sheet and provincias Spain totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------------
# 1. Update the "last updated" timestamp shown in cell A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 15:14"

# ---------------------------------------------------------------------------
# 2. Update per-country statistics (columns B..H) for the rows whose totals
#    changed. Row numbers below match the worksheet row numbers.
# ---------------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6589020
$ws.Range("C4").Value = 857
$ws.Range("D4").Value = 3880707
$ws.Range("E4").Value = 2511968

# Row 19 - Arabia Saudita
$ws.Range("B19").Value = 324407
$ws.Range("C19").Value = 687
$ws.Range("D19").Value = 300933
$ws.Range("E19").Value = 19261
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 4213

# Row 24 - Filipinas
$ws.Range("B24").Value = 258125
$ws.Range("C24").Value = 18
$ws.Range("E24").Value = 15406

# Row 38 - Bolivia
$ws.Range("B38").Value = 93475
$ws.Range("C38").Value = 653
$ws.Range("D38").Value = 83660
$ws.Range("E38").Value = 9258
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 557

# Row 41 - Republica Dominicana
$ws.Range("B41").Value = 86505
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 5846

# Row 44 - now holds "Paises Bajos" data (country order swapped with Marruecos
# because Paises Bajos' total now exceeds Marruecos' total).
$ws.Range("A44").Value = "Paises Bajos"
$ws.Range("B44").Value = 79781
$ws.Range("C44").Value = 1270
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 6252

# Row 45 - now holds "Marruecos" data (its totals are unchanged, it simply
# moved down one row).
$ws.Range("A45").Value = "Marruecos"
$ws.Range("B45").Value = 79767
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 61850
$ws.Range("E45").Value = 16426
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 1491

# Row 47 - Belgica
$ws.Range("B47").Value = 73784
$ws.Range("C47").Value = 193
$ws.Range("D47").Value = 72369
$ws.Range("E47").Value = 677
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 738

# Row 61 - Nigeria
$ws.Range("B61").Value = 46160
$ws.Range("C61").Value = 687
$ws.Range("D61").Value = 42857
$ws.Range("E61").Value = 2926
$ws.Range("G61").Value = 7
$ws.Range("H61").Value = 377

# Row 63 - Nepal
$ws.Range("B63").Value = 45388
$ws.Range("C63").Value = 75
$ws.Range("D63").Value = 44267
$ws.Range("E63").Value = 836
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 285

# Row 67 - Kenia
$ws.Range("B67").Value = 38037
$ws.Range("C67").Value = 163
$ws.Range("D67").Value = 35440
$ws.Range("E67").Value = 2040
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 557

# Row 80
$ws.Range("B80").Value = 21908
$ws.Range("C80").Value = 969
$ws.Range("D80").Value = 2506
$ws.Range("E80").Value = 19050
$ws.Range("G80").Value = 13
$ws.Range("H80").Value = 352

# Row 82
$ws.Range("B82").Value = 19216
$ws.Range("C82").Value = 292
$ws.Range("D82").Value = 16139
$ws.Range("E82").Value = 2448

# Row 149
$ws.Range("B149").Value = 2161
$ws.Range("C149").Value = 4
$ws.Range("D149").Value = 2079
$ws.Range("E149").Value = 72

# Row 194
$ws.Range("B194").Value = 135
$ws.Range("C194").Value = 4
$ws.Range("E194").Value = 79
